# B6-PowerPoint.pptx edit: Mon, Aug 03, 2020 12:05:14 PM
#
# 1) Re-colour the presentation's theme from the custom "Integral" / "Red
#    Violet" palette to the stock Office theme palette.
# 2) Apply the (built-in) table style {0956FA71-CFC9-46D4-9541-6ECDB01A0121}
#    to the three tables that were still using the default "no style, no
#    grid" table style.

$p = $ppt.ActivePresentation

# --- 1) Theme colours -> Office palette ------------------------------------
# ThemeColorScheme indices: 1=dk1 2=lt1 3=dk2 4=lt2 5-10=accent1-6
# 11=hlink 12=folHlink. RGB() isn't available in this host, so the packed
# BGR-ordered integers (r + g*256 + b*65536) are supplied directly.
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72

# --- 2) Table styles ---------------------------------------------------------
$newStyleId = "{0956FA71-CFC9-46D4-9541-6ECDB01A0121}"

$slide14 = $p.Slides.Item(14)
$slide14.Shapes.Item(1).Table.ApplyStyle($newStyleId)

$slide15 = $p.Slides.Item(15)
$slide15.Shapes.Item(1).Table.ApplyStyle($newStyleId)

$slide16 = $p.Slides.Item(16)
$slide16.Shapes.Item(1).Table.ApplyStyle($newStyleId)
